# "Added numbers in remainder report"
# The workbook's single sheet ("Bonus") gets:
#  1. Header row 1 columns C:L rotated left by one (Сим_карты moves from C to L).
#  2. Column A names for rows 2-46 reshuffled to a new order.
#  3. B2 value changed from 14733.5 to 54078.95.
#  4. Five new rows (47-51) appended with names and zeroed B:L data,
#     extending the sheet's used range to A1:L51.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Re-order the category header row (B1 stays put, C1:L1 rotate left) ---
$headers = @("Смартфоны", "Кнопки", "Iphone", "Страховки", "Подписки", "Услуги", "КЭО", "Гаджеты", "Модемы", "Аксы", "Сим_карты")
for ($i = 0; $i -lt $headers.Length; $i++) {
    $col = 2 + $i   # column B = 2 .. column L = 12
    $ws.Cells.Item(1, $col).Value = $headers[$i]
}

# --- 2. Reassign the employee names in column A for the existing rows ---
$renames = @{
    2  = "Цыгина"
    3  = "Чернокрылюк"
    4  = "Лямзина"
    5  = "Морозова"
    7  = "Суворов"
    9  = "Кузякин"
    11 = "Куликова"
    14 = "Романцова"
    16 = "Винокуров"
    17 = "Носова"
    19 = "Буянова"
    20 = "Мелкумян"
    21 = "Боницкий"
    22 = "Светов"
    23 = "Кусочкова"
    24 = "Чекушкин"
    25 = "Рогачев"
    26 = "Илюхина"
    27 = "Караганская"
    28 = "Калашникова"
    29 = "Калугина"
    30 = "Плетникова"
    31 = "Хохлова"
    32 = "Тяпин"
    33 = "Сытин"
    34 = "Голованов"
    35 = "Кузнецов"
    36 = "Винокуров"
    37 = "Михайлова"
    38 = "Сазонов"
    39 = "Александр"
    40 = "Сухарев"
    42 = "Виктория"
    43 = "Миронова"
    44 = "Белозерова"
    45 = "Сироткин"
}
foreach ($row in $renames.Keys) {
    $ws.Cells.Item($row, 1).Value = $renames[$row]
}

# --- 3. Update the bonus total for row 2 (now "Цыгина") ---
$ws.Range("B2").Value = 54078.95

# --- 4. Append the five new trailing rows, all data columns zeroed ---
$newRows = @{
    47 = "Сорманова"
    48 = "Шувалова"
    49 = "Смирнова"
    50 = "Кочетова"
    51 = "Швецова"
}
foreach ($row in $newRows.Keys) {
    $ws.Cells.Item($row, 1).Value = $newRows[$row]
    for ($col = 2; $col -le 12; $col++) {
        $ws.Cells.Item($row, $col).Value = 0
    }
}
